# feat: add 2022-Q4 data
#
# The workbook currently has two sheets: "总计" (summary) and "2022-Q3"
# (fund-holdings detail for 2022-Q3). We need to:
#   1. Keep the existing "2022-Q3" detail data intact (by copying the
#      sheet to a new tab placed right after it), then turn the original
#      tab into the new "2022-Q4" sheet populated with the Q4 figures.
#   2. Insert a new row at the top of the "总计" summary sheet for the
#      2022-Q4 totals, pushing the 2022-Q3 summary row down one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) Sheets: duplicate "2022-Q3" so its data/formatting survive under the
#    same name on a new tab placed right after it, then repurpose the
#    original tab as "2022-Q4" with the new fund-holdings table.
# ---------------------------------------------------------------------
$wsQ3 = $wb.Worksheets.Item("2022-Q3")
$wsQ3.Copy($null, $wsQ3) | Out-Null

$wsQ4 = $wb.Worksheets.Item("2022-Q3")
$wsQ4.Name = "2022-Q4"

$wsQ3New = $wb.Worksheets.Item("2022-Q3 (2)")
$wsQ3New.Name = "2022-Q3"

# Wipe everything on what is now the "2022-Q4" tab so we can lay down the
# fresh Q4 table without any leftover Q3 values/formatting.
$wsQ4.Cells.Clear()

# A donor cell carrying the bold+bordered "header" style used throughout
# this workbook (同"总计" sheet的表头/序号列样式), so the new table picks
# up the exact same look instead of inventing a new one.
$wsTotal = $wb.Worksheets.Item("总计")
$styleDonor = $wsTotal.Range("B1")

# -- header row -----------------------------------------------------
$headers = @{ "B1"="基金代码"; "C1"="基金名称"; "D1"="基金规模"; "E1"="股票总仓位"; "F1"="仓位占比"; "G1"="持有市值(亿元)"; "H1"="仓位排名" }
foreach ($addr in $headers.Keys) {
    $wsQ4.Range($addr).Value = $headers[$addr]
}
$styleDonor.Copy() | Out-Null
$wsQ4.Range("B1:H1").PasteSpecial(-4122) | Out-Null

# -- data rows (text-typed columns B:G keep leading zeros / exact text) --
function Set-TextCell($ws, $addr, $text) {
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $text
    $ws.Range($addr).ClearFormats()
}

Set-TextCell $wsQ4 "B2" "005360"
Set-TextCell $wsQ4 "C2" "汇安资产轮动灵活配置混合A"
Set-TextCell $wsQ4 "D2" "0.26"
Set-TextCell $wsQ4 "E2" "94.03"
Set-TextCell $wsQ4 "F2" "6.04"
Set-TextCell $wsQ4 "G2" "0.0157"
$wsQ4.Range("H2").Value = 10

Set-TextCell $wsQ4 "B3" "017213"
Set-TextCell $wsQ4 "C3" "汇安资产轮动灵活配置混合C"
Set-TextCell $wsQ4 "D3" "0.01"
Set-TextCell $wsQ4 "E3" "94.03"
Set-TextCell $wsQ4 "F3" "6.04"
Set-TextCell $wsQ4 "G3" "0.0006"
$wsQ4.Range("H3").Value = 10

# -- numeric index column (A) with the shared bold+bordered style -------
$wsQ4.Range("A2").Value = 0
$wsQ4.Range("A3").Value = 1
$styleDonor.Copy() | Out-Null
$wsQ4.Range("A2:A3").PasteSpecial(-4122) | Out-Null

# ---------------------------------------------------------------------
# 2) "总计" summary sheet: insert the 2022-Q4 row above the existing
#    2022-Q3 row (which shifts down to row 3 and gets its index bumped).
# ---------------------------------------------------------------------
$wsTotal.Rows.Item(2).Insert()
$wsTotal.Range("B2:D2").ClearFormats()

$wsTotal.Range("A2").Value = 0
$wsTotal.Range("B2").Value = "2022-Q4"
$wsTotal.Range("C2").Value = 2
$wsTotal.Range("D2").Value = 0.02

$wsTotal.Range("A3").Copy() | Out-Null
$wsTotal.Range("A2").PasteSpecial(-4122) | Out-Null
$wsTotal.Range("A3").Value = 1

# Leave the workbook showing the "总计" tab (unchanged from the original),
# same as before this edit.
$wsTotal.Activate()
